# ---------------------------------------------------------------------------
# Applies the "bug fixing, from unique name+bday removed bday" commit:
#   1. On sheet "current": fix a typo'd doctor name, correct a summary count,
#      add three new numeric summary cells, and move the active-cell
#      selection.
#   2. On sheet "2024-06-02": append the 13 patient-visit rows that were
#      captured for that day (previously the sheet only had the header row).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "current"
# ---------------------------------------------------------------------------
$current = $wb.Worksheets.Item("current")

$current.Range("B2").Value = "Karp_Kuzmin_records"
$current.Range("C3").Value = 13
$current.Range("E3").Value = 3
$current.Range("F3").Value = 7
$current.Range("G3").Value = 3

$current.Activate()
$current.Range("B10").Select()

# ---------------------------------------------------------------------------
# Sheet "2024-06-02"
# ---------------------------------------------------------------------------
$day = $wb.Worksheets.Item("2024-06-02")

# Patient visits recorded on 2024-06-02. Column B holds the submission
# timestamp (Excel serial date/time); every other column is free-form text
# captured from a form, so it must stay text even when it looks numeric
# (ids, phone-like "pressure" codes, ISO birth dates, ...).
$rows = @(
  @(2,  "1",  45445.70815164352, "sjfn'pdfi",   "М", "2006-03-16", "opzdkth",       "456879"),
  @(3,  "2",  45445.70842616898, "sjfn'pdfi",   "М", "2006-03-16", "fyckuj",        "435678"),
  @(4,  "3",  45445.71101790509, "aedf",        "М", "2006-05-18", "tjknlsgk",      "435678"),
  @(5,  "4",  45445.7128197338,  "aedf",        "Ж", "2006-06-01", "opzdkth",       "456879"),
  @(6,  "5",  45445.71299002314,"gvkjhkj",      "Р", "2024-05-09", "tjknlsgk",      "456879"),
  @(7,  "6",  45445.79209467593,"jhk",          "Р", "2024-05-31", "fyckuj",        "435678"),
  @(8,  "7",  45445.79244556713,"sjfn'pdfi",    "Р", "2024-05-31", "opzdkth",       "456879"),
  @(9,  "8",  45445.79642175926,"aedf",         "Ж", "2006-05-29", "porrkthso[rhk", "456879"),
  @(10, "9",  45445.81492644676,"sjfn'pdfi",    "Ж", "2006-05-29", "dkjghsropjk",   "456879"),
  @(11, "10", 45445.85901008102,"jhk",          "Ж", "2006-05-30", "dkjghsropjk",   "456879"),
  @(12, "11", 45445.86378913194,"jhk",          "Ж", "2006-05-30", "porrkthso[rhk", "435678"),
  @(13, "12", 45445.8720893287, "jhk",          "Ж", "2006-05-30", "tjknlsgk",      "456879"),
  @(14, "13", 45445.87463234954,"jhkgfguhjk",   "Ж", "2006-05-29", "dkjghsropjk",   "435678")
)

# Timestamps for rows 2-10 share the workbook's pre-existing custom date/time
# format (numFmtId 164); the last four submissions (rows 11-14) get their own
# equivalent format so that style is created fresh (numFmtId 165).
$existingTimeFormat = "yyyy\-mm\-dd\ h:mm:ss"
$newTimeFormat = "yyyy-mm-dd h:mm:ss"

foreach ($r in $rows) {
  $rowIndex = $r[0]

  $idCell = $day.Cells.Item($rowIndex, 1)
  $idCell.Value = "'" + $r[1]

  $timeCell = $day.Cells.Item($rowIndex, 2)
  $timeCell.Value = $r[2]
  if ($rowIndex -le 10) {
    $timeCell.NumberFormat = $existingTimeFormat
  } else {
    $timeCell.NumberFormat = $newTimeFormat
  }

  $day.Cells.Item($rowIndex, 3).Value = $r[3]
  $day.Cells.Item($rowIndex, 4).Value = $r[4]

  $dobCell = $day.Cells.Item($rowIndex, 5)
  $dobCell.Value = "'" + $r[5]

  $day.Cells.Item($rowIndex, 6).Value = $r[6]

  $pressureCell = $day.Cells.Item($rowIndex, 7)
  $pressureCell.Value = "'" + $r[7]
}

Write-Host "edit applied"
